$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.357.78'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.848.58'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '240.43'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').Value = '0.6287'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '0.07610'
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('D9').Value = '0.2919'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '24.63'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').Value = '0.07749'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').Value = '5.022'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').Value = '0.6812'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('D14').Value = '0.00001052'
$ws.Range('E14').Value = '  -4.04%  '
$ws.Range('D15').Value = '83.16'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').Value = '6.132'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '29.385.49'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '229.16'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').Value = '0.9997'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '7.474'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D23').Value = '158.58'
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').Value = '0.1390'
$ws.Range('E24').Value = '  +0.08%  '
$ws.Range('D25').Value = '8.443'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').Value = '17.68'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '1.453'
$ws.Range('E27').Value = '  +10.19%  '
$ws.Range('D28').Value = '1.476'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').Value = '0.05611'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('D30').Value = '4.110'
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('D31').Value = '4.057'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').Value = '1.832'
$ws.Range('E32').Value = '  -1.09%  '
$ws.Range('D33').Value = '1.157'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('D34').Value = '0.7009'
$ws.Range('D35').Value = '2.584'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '1.236.40'
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('D37').Value = '0.01809'
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('D38').Value = '2.728'
$ws.Range('E38').Value = '  -1.94%  '
$ws.Range('D39').Value = '6.424'
$ws.Range('E39').Value = '  -1.46%  '
$ws.Range('D40').Value = '0.9025'
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('D41').Value = '0.9995'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').Value = '101.60'
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('D43').Value = '65.61'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').Value = '7.177'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').Value = '0.3997'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.00000000116'
$ws.Range('E46').Value = '  -2.20%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1157'
$ws.Range('E47').Value = '  +2.16%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.004'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.685'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').Value = '0.05698'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').Value = '0.4627'
$ws.Range('E51').Value = '  -0.15%  '
